# Ankit Rajpoot.xlsx - append two more scraped innings rows (row4 = Dubai/KKR game,
# row5 = Abu Dhabi/Mumbai game) to the bottom of the existing stats table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new numeric-looking values (runs/balls/4s/6s/strike-rate) are
# stored as text, same as the rest of the sheet (numbers are kept as text
# throughout this table).
$ws.Range("A4:K5").NumberFormat = "@"

$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 30 2020"
$ws.Range("C4").Value = "KKR won by 37 runs"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Ankit Rajpoot "
$ws.Range("G4").Value = "7"
$ws.Range("H4").Value = "5"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "140.00"

$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 06 2020"
$ws.Range("C5").Value = "Mumbai won by 57 runs"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Mumbai Indians"
$ws.Range("F5").Value = "Ankit Rajpoot "
$ws.Range("G5").Value = "2"
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "40.00"
